# Commit: Add files via upload
# Appends rows 18-26 (new incident records) to the Optical_Power worksheet
# and extends the sheet dimension accordingly (A1:N17 -> A1:N26).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optical_Power")

$newRows = @(
    @(18, '''5672', '''8/20/2025', 'SOLDADO DE LA INDEPENDENCIA 761', 14, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Cable a baja altura', 1, '{"direccionesNormalizadas": [{"altura": 761, "cod_calle": 20103, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.434330", "y": "-34.567571"}, "direccion": "SOLDADO DE LA INDEPENDENCIA 761, CABA", "nombre_calle": "SOLDADO DE LA INDEPENDENCIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.43433, -34.567571, 'Palermo', 'Capital Sur'),
    @(19, '''5673', '''8/20/2025', 'SOLDADO DE LA INDEPENDENCIA 857', 14, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Cables a baja altura', 1, '{"direccionesNormalizadas": [{"altura": 857, "cod_calle": 20103, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.435519", "y": "-34.566886"}, "direccion": "SOLDADO DE LA INDEPENDENCIA 857, CABA", "nombre_calle": "SOLDADO DE LA INDEPENDENCIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.435519, -34.566886, 'Palermo', 'Capital Sur'),
    @(20, '''5674', '''8/20/2025', 'SOLDADO DE LA INDEPENDENCIA 1057', 14, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Cables a baja altura', 1, '{"direccionesNormalizadas": [{"altura": 1057, "cod_calle": 20103, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.437740", "y": "-34.565577"}, "direccion": "SOLDADO DE LA INDEPENDENCIA 1057, CABA", "nombre_calle": "SOLDADO DE LA INDEPENDENCIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.43774, -34.565577, 'Palermo', 'Capital Sur'),
    @(21, '''5772', '''8/20/2025', 'MALABIA 1648', 14, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Cable sueltos y cortados', 1, '{"direccionesNormalizadas": [{"altura": 1648, "cod_calle": 13016, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.426578", "y": "-34.590461"}, "direccion": "MALABIA 1648, CABA", "nombre_calle": "MALABIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.426578, -34.590461, 'Palermo', 'Capital Sur'),
    @(22, '''5779', '''8/20/2025', 'GUATEMALA 4551', 14, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Cable en panza y cortados', 1, '{"direccionesNormalizadas": [{"altura": 4551, "cod_calle": 7093, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.423098", "y": "-34.587342"}, "direccion": "GUATEMALA 4551, CABA", "nombre_calle": "GUATEMALA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.423098, -34.587342, 'Palermo', 'Capital Sur'),
    @(23, '''5782', '''8/20/2025', 'GUARDIA NACIONAL 1069', 9, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Tendido a baja altura', 1, '{"direccionesNormalizadas": [{"altura": 1069, "cod_calle": 7091, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.493951", "y": "-34.647306"}, "direccion": "GUARDIA NACIONAL 1069, CABA", "nombre_calle": "GUARDIA NACIONAL", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.493951, -34.647306, 'Devoto', 'Capital Norte'),
    @(24, '''5648', '''8/20/2025', 'NAZCA AV. 620', 7, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Tendido a baja altura', 1, '{"direccionesNormalizadas": [{"altura": 620, "cod_calle": 14012, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.472312", "y": "-34.624976"}, "direccion": "NAZCA AV. 620, CABA", "nombre_calle": "NAZCA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.472312, -34.624976, 'Paternal', 'Capital Norte'),
    @(25, '''5646', '''8/20/2025', 'MORON 3036', 7, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Tendido a baja altura', 1, '{"direccionesNormalizadas": [{"altura": 3036, "cod_calle": 13132, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.474687", "y": "-34.624781"}, "direccion": "MORON 3036, CABA", "nombre_calle": "MORON", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.474687, -34.624781, 'Devoto', 'Capital Norte'),
    @(26, '''6927', '''8/20/2025', 'ALTOLAGUIRRE 2520', 12, 'Pendiente ADM', 'Optical Power', 'Pendiente', 'Cables en panza', 1, '{"direccionesNormalizadas": [{"altura": 2520, "cod_calle": 1049, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.491932", "y": "-34.575760"}, "direccion": "ALTOLAGUIRRE 2520, CABA", "nombre_calle": "ALTOLAGUIRRE", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}', -58.491932, -34.57576, 'Paternal', 'Capital Norte')
)

$numericCols = @(4, 9, 11, 12)  # D, I, K, L are numeric (Comuna, Attachments, Coordenada_X, Coordenada_Y)

foreach ($row in $newRows) {
    $r = $row[0]
    for ($col = 1; $col -le 14; $col++) {
        $val = $row[$col]
        if ($numericCols -contains $col) {
            $ws.Cells.Item($r, $col).Value = [double]$val
        } else {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}
